$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data rows (row 2 and row 3) before rewriting the table
$ws.Range("A2:D3").ClearContents()

# New data rows 2-11: Materiaux (A/B) and Isolants (C/D) columns, now
# populated from separate lists rather than paired row-by-row.
$materiaux = @(
    @("Kooltherm Mousse phénolique 0.021 (0.9)_Kingspan", 0.02),
    @("Kooltherm Mousse phénolique 0.02 (0.9)_Kingspan", 0.1),
    @("Kooltherm Mousse phénolique 0.02 (0.9)_Kingspan", 0.12),
    @("Laine minérale", 0.04),
    @("Pare vapeur intérieur", 0.001),
    @("Kooltherm Mousse phénolique 0.021 (0.9)_Kingspan", 0.04),
    @("Laine minérale", 0.07000000000000001),
    @("Myral Therm W22", 0.053),
    @("Kooltherm Mousse phénolique 0.02 (0.9)_Kingspan", 0.06),
    @("Kooltherm Mousse phénolique 0.02 (0.9)_Kingspan", 0.08)
)

for ($i = 0; $i -lt $materiaux.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $materiaux[$i][0]
    $ws.Cells.Item($row, 2).Value = $materiaux[$i][1]
}

# Isolants column only has one entry now (row 2: Bois / 0.012)
$ws.Cells.Item(2, 3).Value = "Bois"
$ws.Cells.Item(2, 4).Value = 0.012
